$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 8.204228666666666
$ws.Range("H2").Value = 24.612686
$ws.Range("I2").Value = 0.07326752815431403
$ws.Range("J2").Value = 0.07326752815431405
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.5290173333333
$ws.Range("N2").Value = 460.587052
$ws.Range("O2").Value = 0.3172206968818489
$ws.Range("P2").Value = 0.317220696881849
$ws.Range("Q2").Value = 1259.587165171297
$ws.Range("R2").Value = 11336.28448654167
$ws.Range("S2").Value = 0.02324197633992198
$ws.Range("T2").Value = 0.02324197633992199

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 8.204228666666666
$ws.Range("H3").Value = 24.612686
$ws.Range("I3").Value = 0.07326752815431403
$ws.Range("J3").Value = 0.07326752815431405
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3487728915577651
$ws.Range("P3").Value = 0.3487728915577651
$ws.Range("Q3").Value = 1384.871359542676
$ws.Range("R3").Value = 12463.84223588409
$ws.Range("S3").Value = 0.02555372765167007
$ws.Range("T3").Value = 0.02555372765167007

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 8.204228666666666
$ws.Range("H4").Value = 24.612686
$ws.Range("I4").Value = 0.07326752815431403
$ws.Range("J4").Value = 0.07326752815431405
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 68.09032333333333
$ws.Range("N4").Value = 204.27097
$ws.Range("O4").Value = 0.1406878008722904
$ws.Range("P4").Value = 0.1406878008722904
$ws.Range("Q4").Value = 558.6285826139355
$ws.Range("R4").Value = 5027.657243525419
$ws.Range("S4").Value = 0.01030784741137906
$ws.Range("T4").Value = 0.01030784741137907

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 8.204228666666666
$ws.Range("H5").Value = 24.612686
$ws.Range("I5").Value = 0.07326752815431403
$ws.Range("J5").Value = 0.07326752815431405
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 93.562673
$ws.Range("N5").Value = 280.688019
$ws.Range("O5").Value = 0.1933186106880956
$ws.Range("P5").Value = 0.1933186106880956
$ws.Range("Q5").Value = 767.6095639565592
$ws.Range("R5").Value = 6908.486075609033
$ws.Range("S5").Value = 0.01416397675134292
$ws.Range("T5").Value = 0.01416397675134292

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 81.515531
$ws.Range("H6").Value = 244.546593
$ws.Range("I6").Value = 0.7279711116319884
$ws.Range("J6").Value = 0.7279711116319885
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 153.5290173333333
$ws.Range("N6").Value = 460.587052
$ws.Range("O6").Value = 0.3172206968818489
$ws.Range("P6").Value = 0.317220696881849
$ws.Range("Q6").Value = 12514.99937183487
$ws.Range("R6").Value = 112634.9943465138
$ws.Range("S6").Value = 0.2309275033417536
$ws.Range("T6").Value = 0.2309275033417536

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 81.515531
$ws.Range("H7").Value = 244.546593
$ws.Range("I7").Value = 0.7279711116319884
$ws.Range("J7").Value = 0.7279711116319885
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 168.7997026666667
$ws.Range("N7").Value = 506.3991080000001
$ws.Range("O7").Value = 0.3487728915577651
$ws.Range("P7").Value = 0.3487728915577651
$ws.Range("Q7").Value = 13759.79739551545
$ws.Range("R7").Value = 123838.1765596391
$ws.Range("S7").Value = 0.2538965895744092
$ws.Range("T7").Value = 0.2538965895744092

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 81.515531
$ws.Range("H8").Value = 244.546593
$ws.Range("I8").Value = 0.7279711116319884
$ws.Range("J8").Value = 0.7279711116319885
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 68.09032333333333
$ws.Range("N8").Value = 204.27097
$ws.Range("O8").Value = 0.1406878008722904
$ws.Range("P8").Value = 0.1406878008722904
$ws.Range("Q8").Value = 5550.418862478356
$ws.Range("R8").Value = 49953.76976230521
$ws.Range("S8").Value = 0.102416654794061
$ws.Range("T8").Value = 0.1024166547940611

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 81.515531
$ws.Range("H9").Value = 244.546593
$ws.Range("I9").Value = 0.7279711116319884
$ws.Range("J9").Value = 0.7279711116319885
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 93.562673
$ws.Range("N9").Value = 280.688019
$ws.Range("O9").Value = 0.1933186106880956
$ws.Range("P9").Value = 0.1933186106880956
$ws.Range("Q9").Value = 7626.810971374363
$ws.Range("R9").Value = 68641.29874236927
$ws.Range("S9").Value = 0.1407303639217646
$ws.Range("T9").Value = 0.1407303639217646

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.363806
$ws.Range("H10").Value = 7.091418
$ws.Range("I10").Value = 0.02110987268797113
$ws.Range("J10").Value = 0.02110987268797113
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 153.5290173333333
$ws.Range("N10").Value = 460.587052
$ws.Range("O10").Value = 0.3172206968818489
$ws.Range("P10").Value = 0.317220696881849
$ws.Range("Q10").Value = 362.9128123466373
$ws.Range("R10").Value = 3266.215311119736
$ws.Range("S10").Value = 0.00669648852516531
$ws.Range("T10").Value = 0.006696488525165313

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.363806
$ws.Range("H11").Value = 7.091418
$ws.Range("I11").Value = 0.02110987268797113
$ws.Range("J11").Value = 0.02110987268797113
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 168.7997026666667
$ws.Range("N11").Value = 506.3991080000001
$ws.Range("O11").Value = 0.3487728915577651
$ws.Range("P11").Value = 0.3487728915577651
$ws.Range("Q11").Value = 399.0097499616827
$ws.Range("R11").Value = 3591.087749655144
$ws.Range("S11").Value = 0.007362551337799982
$ws.Range("T11").Value = 0.007362551337799983

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.363806
$ws.Range("H12").Value = 7.091418
$ws.Range("I12").Value = 0.02110987268797113
$ws.Range("J12").Value = 0.02110987268797113
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 68.09032333333333
$ws.Range("N12").Value = 204.27097
$ws.Range("O12").Value = 0.1406878008722904
$ws.Range("P12").Value = 0.1406878008722904
$ws.Range("Q12").Value = 160.9523148372733
$ws.Range("R12").Value = 1448.57083353546
$ws.Range("S12").Value = 0.002969901565164684
$ws.Range("T12").Value = 0.002969901565164685

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.363806
$ws.Range("H13").Value = 7.091418
$ws.Range("I13").Value = 0.02110987268797113
$ws.Range("J13").Value = 0.02110987268797113
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 93.562673
$ws.Range("N13").Value = 280.688019
$ws.Range("O13").Value = 0.1933186106880956
$ws.Range("P13").Value = 0.1933186106880956
$ws.Range("Q13").Value = 221.164007813438
$ws.Range("R13").Value = 1990.476070320942
$ws.Range("S13").Value = 0.004080931259841154
$ws.Range("T13").Value = 0.004080931259841154

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 19.892761
$ws.Range("H14").Value = 59.67828300000001
$ws.Range("I14").Value = 0.1776514875257265
$ws.Range("J14").Value = 0.1776514875257265
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 153.5290173333333
$ws.Range("N14").Value = 460.587052
$ws.Range("O14").Value = 0.3172206968818489
$ws.Range("P14").Value = 0.317220696881849
$ws.Range("Q14").Value = 3054.116048376858
$ws.Range("R14").Value = 27487.04443539172
$ws.Range("S14").Value = 0.05635472867500804
$ws.Range("T14").Value = 0.05635472867500804

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 19.892761
$ws.Range("H15").Value = 59.67828300000001
$ws.Range("I15").Value = 0.1776514875257265
$ws.Range("J15").Value = 0.1776514875257265
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 168.7997026666667
$ws.Range("N15").Value = 506.3991080000001
$ws.Range("O15").Value = 0.3487728915577651
$ws.Range("P15").Value = 0.3487728915577651
$ws.Range("Q15").Value = 3357.892142019064
$ws.Range("R15").Value = 30221.02927817157
$ws.Range("S15").Value = 0.06196002299388585
$ws.Range("T15").Value = 0.06196002299388585

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 19.892761
$ws.Range("H16").Value = 59.67828300000001
$ws.Range("I16").Value = 0.1776514875257265
$ws.Range("J16").Value = 0.1776514875257265
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 68.09032333333333
$ws.Range("N16").Value = 204.27097
$ws.Range("O16").Value = 0.1406878008722904
$ws.Range("P16").Value = 0.1406878008722904
$ws.Range("Q16").Value = 1354.504528482724
$ws.Range("R16").Value = 12190.54075634451
$ws.Range("S16").Value = 0.02499339710168558
$ws.Range("T16").Value = 0.02499339710168559

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 19.892761
$ws.Range("H17").Value = 59.67828300000001
$ws.Range("I17").Value = 0.1776514875257265
$ws.Range("J17").Value = 0.1776514875257265
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 93.562673
$ws.Range("N17").Value = 280.688019
$ws.Range("O17").Value = 0.1933186106880956
$ws.Range("P17").Value = 0.1933186106880956
$ws.Range("Q17").Value = 1861.219892510153
$ws.Range("R17").Value = 16750.97903259138
$ws.Range("S17").Value = 0.03434333875514699
$ws.Range("T17").Value = 0.03434333875514699

$wb.Save()
